# Weekly price-sheet update: a new weekly record (2022-01-27) is inserted
# at row 240, pushing the existing rows 240-327 down to 241-328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 240 (shifts rows 240:327 down to 241:328).
$ws.Rows(240).Insert()

# Populate the newly inserted row 240 with the new weekly record.
$ws.Range("A240").Value = 10
$ws.Range("B240").Value = "Vega Modelo de Temuco"
$ws.Range("C240").Value = "La Araucanía"
$ws.Range("D240").Value = 44588
$ws.Range("E240").Value = 9
$ws.Range("F240").Value = 100112008
$ws.Range("G240").Value = "Coliflor"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 850
$ws.Range("K240").Value = 1000
$ws.Range("L240").Value = 1000
$ws.Range("M240").Value = 1000
$ws.Range("N240").Value = '$/unidad'
$ws.Range("O240").Value = "Provincia de Cautín"
$ws.Range("P240").Value = 1000
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"
